$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall | Enchanted Silver Ink (item id 27772)
$ws.Range("H28").Value = 1641.5834
$ws.Range("J28").Value = 3000.8
$ws.Range("L28").Value = 3000.8
$ws.Range("N28").Value = -3970.8
# Row 40: Stuck in the Moment | Horn Glue (item id 5505)
$ws.Range("H40").Value = 4144.5
$ws.Range("I40").Value = 2893.4
$ws.Range("J40").Value = 6229.6665
$ws.Range("K40").Value = 2893.4
$ws.Range("L40").Value = 6229.6665
$ws.Range("M40").Value = -2718.4
$ws.Range("N40").Value = -6579.6665
# Row 64: Forged from the Void | Void Glue (item id 5506)
$ws.Range("H64").Value = 5430.727
$ws.Range("J64").Value = 5489.8335
$ws.Range("L64").Value = 5489.8335
$ws.Range("N64").Value = -5985.8335
# Row 67: Dodging the Draft (L) | Void Glue (item id 5506)
$ws.Range("H67").Value = 5430.727
$ws.Range("J67").Value = 5489.8335
$ws.Range("L67").Value = 5489.8335
$ws.Range("N67").Value = -7205.8335
# Row 70: Consecrating Congregation | Holy Water (item id 12604)
$ws.Range("H70").Value = 2023.3846
$ws.Range("J70").Value = 2227.2727
$ws.Range("L70").Value = 6681.8181
$ws.Range("N70").Value = -7221.8181
# Row 73: Curbing the Contagion (L) | Holy Water (item id 12604)
$ws.Range("H73").Value = 2023.3846
$ws.Range("J73").Value = 2227.2727
$ws.Range("L73").Value = 6681.8181
$ws.Range("N73").Value = -8553.8181
# Row 74: Adhesive of Antipathy | Wing Glue (item id 5507)
$ws.Range("H74").Value = 8503.75
$ws.Range("I74").Value = 7505.625
$ws.Range("K74").Value = 7505.625
$ws.Range("M74").Value = -6569.625
# Row 77: It's Gonna Grow Back (L) | Wing Glue (item id 5507)
$ws.Range("H77").Value = 8503.75
$ws.Range("I77").Value = 7505.625
$ws.Range("K77").Value = 37528.125
$ws.Range("M77").Value = -32848.125
# Row 92: Whinier than the Sword | Enchanted Koppranickel Ink (item id 19901)
$ws.Range("H92").Value = 67767.13
$ws.Range("I92").Value = 72379.07000000001
$ws.Range("K92").Value = 72379.07000000001
$ws.Range("M92").Value = -71131.07000000001
# Row 100: Asking for a Friend | Beetle Glue (item id 19906)
$ws.Range("H100").Value = 3638.1333
$ws.Range("I100").Value = 3472.9614
$ws.Range("J100").Value = 4711.75
$ws.Range("K100").Value = 3472.9614
$ws.Range("L100").Value = 4711.75
$ws.Range("M100").Value = -2931.9614
$ws.Range("N100").Value = -5793.75

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot (item id 43999)
$ws.Range("H61").Value = 2841.611
$ws.Range("I61").Value = 2920.5293
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 2920.5293
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -2708.5293
$ws.Range("N61").Value = -1924
# Row 136: Metal with Mettle | Cobalt Tungsten Ingot (item id 43999)
$ws.Range("H136").Value = 2841.611
$ws.Range("I136").Value = 2920.5293
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 8761.5879
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -6211.5879
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt | Iron Ingot (item id 14149)
$ws.Range("H20").Value = 5884904.5
$ws.Range("I20").Value = 12501318
$ws.Range("J20").Value = 3648.889
$ws.Range("K20").Value = 12501318
$ws.Range("L20").Value = 3648.889
$ws.Range("M20").Value = -12501071
$ws.Range("N20").Value = -4142.889
# Row 86: Through Thick and Thin | Adamantite Nugget (item id 12526)
$ws.Range("H86").Value = 4085.6667
$ws.Range("I86").Value = 1966
$ws.Range("J86").Value = 4509.6
$ws.Range("K86").Value = 1966
$ws.Range("L86").Value = 4509.6
$ws.Range("M86").Value = -843
$ws.Range("N86").Value = -6755.6
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget (item id 12526)
$ws.Range("H89").Value = 4085.6667
$ws.Range("I89").Value = 1966
$ws.Range("J89").Value = 4509.6
$ws.Range("K89").Value = 9830
$ws.Range("L89").Value = 22548
$ws.Range("M89").Value = -4214
$ws.Range("N89").Value = -33780

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber (item id 44023)
$ws.Range("H31").Value = 1082.8462
$ws.Range("I31").Value = 1055.4
$ws.Range("J31").Value = 1100
$ws.Range("K31").Value = 1055.4
$ws.Range("L31").Value = 1100
$ws.Range("M31").Value = -760.4000000000001
$ws.Range("N31").Value = -1690
# Row 34: Armoires of the Rich and Famous | Walnut Lumber (item id 44023)
$ws.Range("H34").Value = 1082.8462
$ws.Range("I34").Value = 1055.4
$ws.Range("J34").Value = 1100
$ws.Range("K34").Value = 1055.4
$ws.Range("L34").Value = 1100
$ws.Range("M34").Value = -853.4000000000001
$ws.Range("N34").Value = -1504
# Row 62: Splinter in the Sewers | Cedar Lumber (item id 12580)
$ws.Range("H62").Value = 4020.8
$ws.Range("I62").Value = 4699.3335
$ws.Range("K62").Value = 4699.3335
$ws.Range("M62").Value = -4075.3335
# Row 65: The Lumber of Their Discontent (L) | Cedar Lumber (item id 12580)
$ws.Range("H65").Value = 4020.8
$ws.Range("I65").Value = 4699.3335
$ws.Range("K65").Value = 23496.6675
$ws.Range("M65").Value = -20376.6675

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food | Table Salt (item id 4847)
$ws.Range("H2").Value = 49.615383
$ws.Range("I2").Value = 44.375
$ws.Range("K2").Value = 266.25
$ws.Range("M2").Value = -153.25
# Row 33: Cooking with Gas | Chicken Stock (item id 4867)
$ws.Range("H33").Value = 70
$ws.Range("I33").Value = 70
$ws.Range("K33").Value = 420
$ws.Range("M33").Value = -137
# Row 40: True Grits | Cornmeal (item id 4827)
$ws.Range("H40").Value = 168.875
$ws.Range("J40").Value = 211.83333
$ws.Range("L40").Value = 847.33332
$ws.Range("N40").Value = -985.33332
# Row 125: At Any Temperature | Borscht (item id 36043)
$ws.Range("H125").Value = 6750
$ws.Range("I125").Value = 6750
$ws.Range("K125").Value = 20250
$ws.Range("M125").Value = -15330

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle | Durium Ingot (item id 36169)
$ws.Range("H102").Value = 4252.1333
$ws.Range("I102").Value = 2997.8572
$ws.Range("J102").Value = 5349.625
$ws.Range("K102").Value = 2997.8572
$ws.Range("L102").Value = 5349.625
$ws.Range("M102").Value = -1375.8572
$ws.Range("N102").Value = -8593.625
# Row 132: On Board for Lar | Lar Ingot (item id 44008)
$ws.Range("H132").Value = 169417.83
$ws.Range("I132").Value = 202701.6
$ws.Range("K132").Value = 608104.8
$ws.Range("M132").Value = -605574.8

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad | Toad Leather (item id 36248)
$ws.Range("H40").Value = 8565.429
$ws.Range("I40").Value = 6369.5
$ws.Range("K40").Value = 6369.5
$ws.Range("M40").Value = -6233.5
# Row 132: Tenets of Tanning | Silver Lobo Leather (item id 44058)
$ws.Range("H132").Value = 154481
$ws.Range("I132").Value = 202474.83
$ws.Range("J132").Value = 10499.5
$ws.Range("K132").Value = 607424.49
$ws.Range("L132").Value = 31498.5
$ws.Range("M132").Value = -604894.49
$ws.Range("N132").Value = -36558.5

$ws = $wb.Worksheets.Item("WVR")
# Row 21: Don't Trew So Hard | Initiate's Slops (item id 3341)
$ws.Range("H21").Value = 7994.6665
$ws.Range("J21").Value = 7994.6665
$ws.Range("L21").Value = 7994.6665
$ws.Range("N21").Value = -8464.666499999999
# Row 35: Pantser Corps | Initiate's Slops (item id 3341)
$ws.Range("H35").Value = 7994.6665
$ws.Range("J35").Value = 7994.6665
$ws.Range("L35").Value = 7994.6665
$ws.Range("N35").Value = -8574.666499999999
# Row 126: A Polished Purchase | Snow Linen (item id 36210)
$ws.Range("H126").Value = 131431.75
$ws.Range("I126").Value = 148349.86
$ws.Range("J126").Value = 13005
$ws.Range("K126").Value = 445049.58
$ws.Range("L126").Value = 39015
$ws.Range("M126").Value = -442579.58
$ws.Range("N126").Value = -43955
# Row 132: Comfy Cabins | Snow Cotton Cloth (item id 44029)
$ws.Range("H132").Value = 351999.34
$ws.Range("I132").Value = 502999.5
$ws.Range("K132").Value = 1508998.5
$ws.Range("M132").Value = -1506468.5
